# The author clicked in the middle of the final sentence ("...agile
# method. We will also be using integration for some  modules.") and
# retyped the duplicate space out, which is exactly the kind of edit
# that makes Word split a run into pieces around the insertion point
# and relocate the "_GoBack" bookmark Word maintains for the last
# edited location.
$d = $word.ActiveDocument

# Locate the run that holds the sentence we need to touch.
$found = $d.Content
$found.Find.Execute("od. We will also be using integration for some  modules.",
                     $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runStart = $found.Start

# Offsets (relative to $runStart) of the two edit points inside that run.
$splitAfterAls = "od. We will als".Length
$beforeDoubleSpace = "od. We will also be using integration for some".Length

$splitPos = $runStart + $splitAfterAls
$spacePos = $runStart + $beforeDoubleSpace

# 1) Split "od. We will als" | "o be using integration for some  modules."
#    by dropping a temporary bookmark at the boundary -- a bookmark start/end
#    pair forces Word to end one run and begin another at that character
#    position. Leave the marker in place for now so this split survives the
#    edits that follow.
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TempSplitMarker", $splitRange)

# 2) Remove one of the two spaces before "modules." (the actual text edit).
$spaceRange = $d.Range($spacePos, $spacePos + 1)
$spaceRange.Text = ""

# 3) Word keeps a single "_GoBack" bookmark marking the most recent edit
#    location; re-adding it moves it here (right before " modules.") and
#    splits the run at that point too.
$goBack = $d.Range($spacePos, $spacePos)
$d.Bookmarks.Add("_GoBack", $goBack)

# 4) Now that both run boundaries exist, drop the temporary marker.
$d.Bookmarks("TempSplitMarker").Delete()
